# Materialliste.xlsx update:
#   - add two new rows (M12 Gewindestange / M12 Stopmuttern) to the
#     "ebay - Sonstige Kaeufe" table on Tabelle2
#   - move the subtotal/grand-total block down to make room
#   - refresh the two summary formulas

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

$CUR = '_-* #,##0.00" €"_-;-* #,##0.00" €"_-;_-* "-"?? " €"_-;_-@_-'
$EUR = '#,##0.00" €";[Red]-#,##0.00" €"'

# ---------------------------------------------------------------
# New row 114: M12 Gewindestange
# ---------------------------------------------------------------
$ws.Range("A114").Value = 3
$ws.Range("B114").Value = "M12 Gewindestange"
$ws.Range("C114").Value = "M12 1m"
$ws.Range("D114").Value = 3.13
$ws.Range("D114").NumberFormat = $CUR
$ws.Range("E114").Formula = "=A114*D114"
$ws.Range("E114").NumberFormat = $EUR
$ws.Range("F114").Value = "https://www.ebay.de/itm/Gewindestange-A2-V2A-Niro-Edelstahl-DIN-975-976-Gewindestangen-1000-mm-1m-Meter-/271358495006?var=&hash=item3f2e3a091e:m:mywm3Vn9Ce40FjrElud3qbw"
$ws.Range("G114").Value = "Verstärkung der Z-Achse"

# ---------------------------------------------------------------
# New row 115: M12 Stopmuttern
# ---------------------------------------------------------------
$ws.Range("A115").Value = 25
$ws.Range("B115").Value = "M12 Stopmuttern"
$ws.Range("C115").Value = "M12 1m"
$ws.Range("D115").Value = 9.4499999999999993
$ws.Range("D115").NumberFormat = $CUR
$ws.Range("E115").Value = 9.4499999999999993
$ws.Range("E115").NumberFormat = $CUR
$ws.Range("F115").Value = "https://www.ebay.de/itm/Stopmutter-Sicherungsmuttern-DIN985-Edelstahl-VA-M2-M3-M4-M5-M6-M8-M12-Polystop/301329119935?epid=11029146985&hash=item46289d7abf:m:mwdeZqra2fNfOcTC50amyKA"

# ---------------------------------------------------------------
# Clear the old row 117 label/grand-total, keep "zzgl. Versand" in F117
# and put the section subtotal (previously on row 115) on row 117.
# ---------------------------------------------------------------
$ws.Range("D117").ClearContents()
$ws.Range("E117").Formula = "=SUM(E88:E115)"
$ws.Range("E117").Style = "Ausgabe"
$ws.Range("E117").NumberFormat = $CUR

# ---------------------------------------------------------------
# Push the two trailing notes down from rows 119/120 to 121/122
# (read before row 119 is overwritten with the grand-total label)
# ---------------------------------------------------------------
$ws.Range("D121").Value = $ws.Range("D119").Value()
$ws.Range("D122").Value = $ws.Range("D120").Value()
$ws.Range("D120").ClearContents()

# ---------------------------------------------------------------
# Row 119: grand total (previously on row 117), with updated reference
# ---------------------------------------------------------------
$ws.Range("D119").Value = "Gesamter Betrag: "
$ws.Range("D119").Style = "Ausgabe"
$ws.Range("E119").Formula = "=E30+E40+E51+E62+E73+E83+E117"
$ws.Range("E119").Style = "Ausgabe"
$ws.Range("E119").NumberFormat = $EUR

# ---------------------------------------------------------------
# View state
# ---------------------------------------------------------------
$ws.Range("D114").Select()
$ws.Application.ActiveWindow.ScrollRow = 97
